$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values as plain text in this sheet (e.g. "24.446.82").
# Temporarily force the range to Text format so numeric-looking updates
# ("1.001", "1.000", "0.9997", ...) are not auto-converted to numbers,
# then restore the original (Normal / General) style once all values are set.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '24.432.94'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '1.670.35'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '312.67'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').Value = '0.3963'
$ws.Range('E7').Value = '  +1.51%  '
$ws.Range('D8').Value = '0.3932'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('D9').Value = '52.65'
$ws.Range('E9').Value = '  +6.94%  '
$ws.Range('D10').Value = '1.397'
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('D11').Value = '1.000'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').Value = '0.08571'
$ws.Range('E12').Value = '  +2.21%  '
$ws.Range('D13').Value = '24.58'
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').Value = '7.313'
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('D15').Value = '7.950'
$ws.Range('E15').Value = '  +7.00%  '
$ws.Range('D16').Value = '0.00001338'
$ws.Range('E16').Value = '  +5.25%  '
$ws.Range('D17').Value = '1.664.54'
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '95.15'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '0.07033'
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').Value = '20.65'
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('D21').Value = '7.002'
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '13.78'
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('D24').Value = '24.431.63'
$ws.Range('E24').Value = '  +1.62%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.462'
$ws.Range('E25').Value = '  +6.09%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '3.067'
$ws.Range('E26').Value = '  +14.34%  '
$ws.Range('D27').Value = '22.57'
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('D28').Value = '157.59'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = '143.00'
$ws.Range('E29').Value = '  +2.30%  '
$ws.Range('D30').Value = '5.446'
$ws.Range('E30').Value = '  +2.86%  '
$ws.Range('D31').Value = '8.040'
$ws.Range('E31').Value = '  -6.55%  '
$ws.Range('D32').Value = '2.541'
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('D33').Value = '1.843.85'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = '1.064'
$ws.Range('E34').Value = '  +12.91%  '
$ws.Range('D35').Value = '0.03109'
$ws.Range('E35').Value = '  +7.78%  '
$ws.Range('D36').Value = '0.08291'
$ws.Range('E36').Value = '  +4.43%  '
$ws.Range('D37').Value = '6.923'
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('D38').Value = '11.22'
$ws.Range('E38').Value = '  +13.83%  '
$ws.Range('D39').Value = '0.2771'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').Value = '0.09276'
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('D41').Value = '0.7730'
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('D42').Value = '13.83'
$ws.Range('E42').Value = '  +6.76%  '
$ws.Range('D43').Value = '1.446'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '16.69'
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('D45').Value = '0.7129'
$ws.Range('E45').Value = '  +4.29%  '
$ws.Range('D46').Value = '2.555'
$ws.Range('E46').Value = '  +4.28%  '
$ws.Range('D47').Value = '4.131'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '0.9994'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('D49').Value = '0.08455'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').Value = '137.13'
$ws.Range('E50').Value = '  +4.40%  '
$ws.Range('D51').Value = '1.272'
$ws.Range('E51').Value = '  +2.41%  '

$priceRange.Style = "Normal"
